# Applies the "Updated cryptos list" data refresh described in the commit.
# Rows 29/30, 31/32 and 41/42 were re-ranked (their Coin/Link pairs swapped),
# and every row's Price (D) and Volume(1h) (E) were refreshed with new quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.747.74"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "2.483.95"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'575.79"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "'149.43"
$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.540"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").Value = "2.484.36"
$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").Value = "'5.28"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("D14").Value = "'27.23"
$ws.Range("E14").Value = "  +1.28%  "

$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "2.943.69"
$ws.Range("E16").Value = "  +3.87%  "

$ws.Range("D17").Value = "63.594.76"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").Value = "2.473.20"
$ws.Range("E18").Value = "  +1.73%  "

$ws.Range("D19").Value = "'11.55"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  +5.33%  "

$ws.Range("D21").Value = "'330.67"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("D22").Value = "'4.21"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").Value = "'2.09"
$ws.Range("E23").Value = "  +19.64%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'66.22"
$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("D26").Value = "'626.30"
$ws.Range("E26").Value = "  +11.59%  "

$ws.Range("D27").Value = "'0.0000105"
$ws.Range("E27").Value = "  +5.32%  "

$ws.Range("D28").Value = "'8.65"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.664.07"
$ws.Range("E29").Value = "  +4.24%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.54"
$ws.Range("E30").Value = "  +5.01%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.41"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").Value = "'0.143"
$ws.Range("E33").Value = "  -3.14%  "

$ws.Range("D34").Value = "'1.91"
$ws.Range("E34").Value = "  +1.88%  "

$ws.Range("D35").Value = "'5.25"
$ws.Range("E35").Value = "  +7.27%  "

$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'0.385"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").Value = "'5.49"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").Value = "'18.84"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'148.49"
$ws.Range("E41").Value = "  -1.23%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.83"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").Value = "'2.72"
$ws.Range("E43").Value = "  +12.64%  "

$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").Value = "'149.86"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "'3.77"
$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("D47").Value = "'21.37"
$ws.Range("E47").Value = "  +4.42%  "

$ws.Range("D48").Value = "'0.0546"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("D49").Value = "'0.608"
$ws.Range("E49").Value = "  +0.96%  "

$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  -0.86%  "
